$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2100
$ws.Range("J31").Value = 2900
$ws.Range("L31").Value = 8700
$ws.Range("N31").Value = -9160

$ws.Range("H41").Value = 131.9
$ws.Range("I41").Value = 128.42857
$ws.Range("J41").Value = 140
$ws.Range("K41").Value = 128.42857
$ws.Range("L41").Value = 140
$ws.Range("M41").Value = 311.57143
$ws.Range("N41").Value = -1020

$ws.Range("H62").Value = 2400
$ws.Range("I62").Value = 1550
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 1550
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -926
$ws.Range("N62").Value = -4498

$ws.Range("H65").Value = 2400
$ws.Range("I65").Value = 1550
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 7750
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -4630
$ws.Range("N65").Value = -22490

$ws.Range("H112").Value = 3833.182
$ws.Range("J112").Value = 4011.5
$ws.Range("L112").Value = 12034.5
$ws.Range("N112").Value = -14250.5

$ws.Range("H138").Value = 4941.122
$ws.Range("I138").Value = 6399.5713
$ws.Range("J138").Value = 4640.853
$ws.Range("K138").Value = 19198.7139
$ws.Range("L138").Value = 13922.559
$ws.Range("M138").Value = -14058.7139
$ws.Range("N138").Value = -24202.559

$ws.Range("H141").Value = 4596.4546
$ws.Range("I141").Value = 2385.923
$ws.Range("J141").Value = 7789.4443
$ws.Range("K141").Value = 7157.768999999999
$ws.Range("L141").Value = 23368.3329
$ws.Range("M141").Value = -1977.768999999999
$ws.Range("N141").Value = -33728.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 519812.47
$ws.Range("I32").Value = 680098.4
$ws.Range("J32").Value = 14295.385
$ws.Range("K32").Value = 680098.4
$ws.Range("L32").Value = 14295.385
$ws.Range("M32").Value = -679811.4
$ws.Range("N32").Value = -14869.385

$ws.Range("H97").Value = 1169.8125
$ws.Range("I97").Value = 1095.238
$ws.Range("J97").Value = 1312.1818
$ws.Range("K97").Value = 1095.238
$ws.Range("L97").Value = 1312.1818
$ws.Range("M97").Value = -599.2380000000001
$ws.Range("N97").Value = -2304.1818

$ws.Range("H122").Value = 1555.3
$ws.Range("I122").Value = 1256
$ws.Range("J122").Value = 1630.125
$ws.Range("K122").Value = 3768
$ws.Range("L122").Value = 4890.375
$ws.Range("M122").Value = -1318
$ws.Range("N122").Value = -9790.375

$ws.Range("H132").Value = 3554.796
$ws.Range("I132").Value = 2734.1316
$ws.Range("J132").Value = 6389.8184
$ws.Range("K132").Value = 8202.3948
$ws.Range("L132").Value = 19169.4552
$ws.Range("M132").Value = -5672.3948
$ws.Range("N132").Value = -24229.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1593.1333
$ws.Range("I31").Value = 1732
$ws.Range("K31").Value = 1732
$ws.Range("M31").Value = -1437

$ws.Range("H34").Value = 1593.1333
$ws.Range("I34").Value = 1732
$ws.Range("K34").Value = 1732
$ws.Range("M34").Value = -1530

$ws.Range("H62").Value = 3801.4285
$ws.Range("I62").Value = 3935
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3935
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -3311
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 3801.4285
$ws.Range("I65").Value = 3935
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 19675
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -16555
$ws.Range("N65").Value = -21240

$ws.Range("H134").Value = 1036.2142
$ws.Range("I134").Value = 750.5833
$ws.Range("K134").Value = 2251.7499
$ws.Range("M134").Value = 283.2501000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 838.4138
$ws.Range("J68").Value = 964.2105
$ws.Range("L68").Value = 2892.6315
$ws.Range("N68").Value = -4514.6315

$ws.Range("H71").Value = 838.4138
$ws.Range("J71").Value = 964.2105
$ws.Range("L71").Value = 8677.8945
$ws.Range("N71").Value = -16789.8945

$ws.Range("H113").Value = 2619.8
$ws.Range("I113").Value = 833
$ws.Range("J113").Value = 5300
$ws.Range("K113").Value = 2499
$ws.Range("L113").Value = 15900
$ws.Range("M113").Value = -329
$ws.Range("N113").Value = -20240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 20142858
$ws.Range("I14").Value = 20142858
$ws.Range("K14").Value = 20142858
$ws.Range("M14").Value = -20142690

$ws.Range("H97").Value = 2950
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2950
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2950
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3942

$ws.Range("H122").Value = 2904.3333
$ws.Range("I122").Value = 1943.909
$ws.Range("J122").Value = 3960.8
$ws.Range("K122").Value = 5831.727000000001
$ws.Range("L122").Value = 11882.4
$ws.Range("M122").Value = -3381.727000000001
$ws.Range("N122").Value = -16782.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5833.8887
$ws.Range("I61").Value = 6840
$ws.Range("J61").Value = 4576.25
$ws.Range("K61").Value = 6840
$ws.Range("L61").Value = 4576.25
$ws.Range("M61").Value = -6638
$ws.Range("N61").Value = -4980.25

$ws.Range("H100").Value = 3168.9285
$ws.Range("I100").Value = 2640.4285
$ws.Range("J100").Value = 3697.4285
$ws.Range("K100").Value = 2640.4285
$ws.Range("L100").Value = 3697.4285
$ws.Range("M100").Value = -2099.4285
$ws.Range("N100").Value = -4779.4285

$ws.Range("H113").Value = 5833.8887
$ws.Range("I113").Value = 6840
$ws.Range("J113").Value = 4576.25
$ws.Range("K113").Value = 6840
$ws.Range("L113").Value = 4576.25
$ws.Range("M113").Value = -4670
$ws.Range("N113").Value = -8916.25

$ws.Range("H136").Value = 1749.4
$ws.Range("I136").Value = 1832.6666
$ws.Range("K136").Value = 5497.9998
$ws.Range("M136").Value = -2947.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1376.0869
$ws.Range("I113").Value = 1460.091
$ws.Range("K113").Value = 4380.272999999999
$ws.Range("M113").Value = -2210.272999999999
